$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. These are text-typed cells in the
# source workbook (prices / percentage-change strings), so we force the
# cell to Text format before writing, then restore the default "Normal"
# style so no stray number-format style is left behind.
$updates = @{
    'D2' = '289.80'
    'E2' = '-3.87%'
    'D3' = '30.86'
    'E3' = '-4.04%'
    'D4' = '4.882'
    'E4' = '-2.31%'
    'D5' = '0.07155'
    'E5' = '-9.36%'
    'D6' = '1.842'
    'E6' = '-12.52%'
    'D7' = '7.638'
    'E7' = '-2.13%'
    'D8' = '3.768'
    'E8' = '-1.54%'
    'E9' = '-3.54%'
    'D10' = '0.1640'
    'E10' = '-6.40%'
    'D11' = '0.07534'
    'E11' = '-5.60%'
    'D12' = '0.08144'
    'E12' = '-5.46%'
    'D13' = '0.02990'
    'E13' = '-3.88%'
    'D14' = '0.09990'
    'E14' = '-0.25%'
    'D15' = '0.001495'
    'E15' = '-2.03%'
    'D16' = '0.005770'
    'E16' = '-3.94%'
    'E18' = '-0.16%'
    'D19' = '2.112'
    'E19' = '-7.20%'
    'D20' = '0.3277'
    'E20' = '-0.33%'
    'E21' = '-1.35%'
    'D22' = '4.271'
    'E22' = '-0.33%'
    'D23' = '0.2001'
    'E23' = '11.65%'
    'D24' = '0.04479'
    'E24' = '-2.60%'
    'D25' = '0.001213'
    'E25' = '-2.04%'
    'D26' = '0.004659'
    'E26' = '4.72%'
    'E27' = '0.07%'
    'D39' = '0.01637'
    'E39' = '-4.64%'
    'E40' = '-9.03%'
    'D41' = '0.007393'
    'E41' = '-1.27%'
    'D42' = '0.1306'
    'E42' = '-3.83%'
    'D43' = '0.002005'
    'E43' = '-12.09%'
    'D44' = '0.01023'
    'E44' = '-0.33%'
    'D45' = '0.00005861'
    'E45' = '-2.26%'
    'E46' = '-0.02%'
    'D47' = '2.207'
    'E47' = '169.06%'
    'E48' = '-11.55%'
    'E49' = '-0.02%'
    'E50' = '-0.02%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
